{"js": "const body = context.document.body;\n\n// The \"_GoBack\" bookmark currently sits at the end of the \"3rd level\" bullet\n// paragraph; it needs to move to the Collectibles paragraph, right after the\n// newly-inserted \"diamond shaped pieces\" text. Remove the old one first since\n// bookmark names must stay unique.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Replace \"circular wisps\" with \"diamond shaped pieces\" in the Collectibles bullet.\nconst target = body.search(\"circular wisps\", { matchCase: true, matchWholeWord: false });\ntarget.load(\"text\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  target.items[0].insertText(\"diamond shaped pieces\", \"Replace\");\n  await context.sync();\n}\n\n// Re-locate the replaced phrase so we can split it into its own run (matching\n// how Word itself leaves a dedicated run around freshly typed text) and plant\n// the relocated \"_GoBack\" bookmark right after it.\nconst replaced = body.search(\"diamond shaped pieces\", { matchCase: true, matchWholeWord: false });\nreplaced.load(\"text\");\nawait context.sync();\n\nconst phrase = replaced.items[0];\n// A transient bookmark forces a run break before the phrase; delete just the\n// bookmark afterwards so only the run split remains.\nphrase.getRange(\"Before\").insertBookmark(\"_TmpRunSplit\");\nphrase.getRange(\"After\").insertBookmark(\"_GoBack\");\nawait context.sync();\n\ncontext.document.deleteBookmark(\"_TmpRunSplit\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The \"_GoBack\" bookmark currently sits at the end of the \"3rd level\" bullet\n# paragraph; it needs to move to the Collectibles paragraph, right after the\n# newly-inserted \"diamond shaped pieces\" text. Remove the old one first since\n# bookmark names must stay unique.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Replace \"circular wisps\" with \"diamond shaped pieces\" in the Collectibles bullet.\n$findRange = $d.Content\n$findRange.Find.Execute(\"circular wisps\", $false, $false, $false, $false, $false, $true, 1, $false, \"diamond shaped pieces\", 2)\n\n# Re-locate the replaced phrase so it can be split into its own run (matching\n# how Word leaves a dedicated run around freshly typed text) and plant the\n# relocated \"_GoBack\" bookmark right after it.\n$phraseRange = $d.Content\n$phraseRange.Find.Execute(\"diamond shaped pieces\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n\n$startPoint = $d.Range($phraseRange.Start, $phraseRange.Start)\n$endPoint = $d.Range($phraseRange.End, $phraseRange.End)\n\n# A transient bookmark forces a run break before the phrase; delete just the\n# bookmark afterwards so only the run split remains.\n$d.Bookmarks.Add(\"_TmpRunSplit\", $startPoint)\n$d.Bookmarks.Add(\"_GoBack\", $endPoint)\n\n$d.Bookmarks(\"_TmpRunSplit\").Delete()\n"}
